# The sheet had two "spacer" columns (C and F) that only carried an
# empty/blank-style cell in every data row and no data of their own.
# Cleaning up the verification-data sheet removes those two empty
# columns entirely, so the remaining data columns (which used to be
# A,B,D,E,G,H) collapse left into A,B,C,D,E,F.
#
# Deleting the higher-numbered column first keeps column 3 (C) valid
# for the second delete.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(6).Delete() | Out-Null   # old column F (empty spacer)
$ws.Columns.Item(3).Delete() | Out-Null   # old column C (empty spacer)

# Reflect the author's final selection/cursor position on the sheet.
$ws.Range("H10").Select() | Out-Null
